# S18/G04: DSL-backed indicator alerts end-to-end
# - Updates the text of the two existing Alerts-page remark/pending-work cells
#   (S18_G03_TF002 row) to reflect the new edit/pause/delete actions.
# - Appends a new "S18 / G04" group block (4 task rows) describing the
#   DSL-backed indicator alerts work, mirroring the existing row layout
#   (sprint#, group#, group task description, task#, task description,
#   deviations, status, remarks, pending work).
# - Nudges a few existing row heights (cosmetic autofit drift) and updates
#   the sheet's scroll/selection position to the new bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the existing S18_G03_TF002 "remarks" (col H) / "pending work"
#    (col I) cells to describe the new edit/pause/delete actions on the
#    Alerts page. Find the row by scanning column D (task#) for the id
#    rather than hard-coding a row number.
# ---------------------------------------------------------------------

$targetRow = 0
$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
for ($r = 1; $r -le $rowCount; $r++) {
    $taskId = $ws.Cells.Item($r, 4).Text
    if ($taskId -eq "S18_G03_TF002") {
        $targetRow = $r
    }
}

if ($targetRow -gt 0) {
    $ws.Cells.Item($targetRow, 8).Value = "Alerts page lists indicator rules with strategy linkage, status, timestamps, and now includes edit/pause/delete management actions for internal alerts."
    $ws.Cells.Item($targetRow, 9).Value = "Integrate TradingView alerts into the Alerts page and add richer filters and bulk actions."
}

# ---------------------------------------------------------------------
# 2. Append the new S18 / G04 group - four task rows (TB001, TB002,
#    TF001, TF002) right after the last existing row.
# ---------------------------------------------------------------------

$lastRow = $usedRange.Rows.Count

$newRows = @(
    @{
        A = "S18"; B = "G04"; C = "DSL-backed indicator alerts end-to-end";
        D = "S18_G04_TB001";
        E = "Extend indicator alert API/model so rules can optionally carry a dsl_expression that is parsed into expression_json while keeping conditions_json backward compatible.";
        F = "Builds directly on the S18/G01 expression engine and keeps existing single-condition alerts working unchanged.";
        G = "implemented";
        H = "indicator_rules now support optional dsl_expression; create/update APIs parse DSL into expression_json while keeping conditions_json for backwards compatibility.";
        I = "Consider allowing rules that rely solely on expression_json without placeholder conditions.";
        Height = 55.2
    },
    @{
        A = "S18"; B = "G04"; C = "DSL-backed indicator alerts end-to-end";
        D = "S18_G04_TB002";
        E = "Wire expression_json evaluation into the indicator alert scheduler so rules backed by the AST drive firing logic, with a safe fallback to legacy condition-based evaluation.";
        F = "Scheduler chooses expression_json when present and falls back to existing conditions_json logic to avoid breaking older rules.";
        G = "implemented";
        H = "Indicator alert scheduler evaluates expression_json-backed rules via the AST engine and falls back to legacy per-condition logic when no expression is present, updating last_evaluated_at for each rule.";
        I = "Tune scheduling/logging and consolidate metrics for expression vs legacy rules.";
        Height = 41.75
    },
    @{
        A = "S18"; B = "G04"; C = "DSL-backed indicator alerts end-to-end";
        D = "S18_G04_TF001";
        E = "Add Simple/DSL mode toggle and DSL editor to the indicator alert dialog, with inline validation and parse errors surfaced from the backend.";
        F = "UI will keep the current single-condition builder as the Simple mode and introduce a DSL textarea tab that talks to a small /api/indicator-alerts/parse endpoint.";
        G = "implemented";
        H = "Holdings alert dialog now has Simple/DSL tabs; DSL mode sends dsl_expression to the backend and uses the same trigger/action settings as the simple builder.";
        I = "Optional: add live DSL validation/preview endpoint and nicer formatting helpers.";
        Height = 41.75
    },
    @{
        A = "S18"; B = "G04"; C = "DSL-backed indicator alerts end-to-end";
        D = "S18_G04_TF002";
        E = "Load and display existing DSL-backed rules in Holdings and Alerts views, including showing the compiled DSL text and allowing edits via the DSL editor.";
        F = "Focus first on read/edit for per-symbol indicator rules; group/basket DSL support can follow later.";
        G = "implemented";
        H = "DSL-backed rules are surfaced in Holdings (Existing alerts list) and on the Alerts page, which now supports editing enabled/trigger mode and updating DSL text, plus delete actions.";
        I = "Expose DSL snippets directly in the Alerts grid and add deeper filters (e.g., by strategy or DSL presence).";
        Height = 41.75
    }
)

$row = $lastRow
foreach ($rowData in $newRows) {
    $row = $row + 1
    $ws.Cells.Item($row, 1).Value = $rowData.A
    $ws.Cells.Item($row, 2).Value = $rowData.B
    $ws.Cells.Item($row, 3).Value = $rowData.C
    $ws.Cells.Item($row, 4).Value = $rowData.D
    $ws.Cells.Item($row, 5).Value = $rowData.E
    $ws.Cells.Item($row, 6).Value = $rowData.F
    $ws.Cells.Item($row, 7).Value = $rowData.G
    $ws.Cells.Item($row, 8).Value = $rowData.H
    $ws.Cells.Item($row, 9).Value = $rowData.I
    $ws.Rows.Item($row).RowHeight = $rowData.Height
}

# ---------------------------------------------------------------------
# 3. Small cosmetic row-height adjustments on the rows just above the new
#    block (autofit drift observed after the edit in the source workbook).
# ---------------------------------------------------------------------

$ws.Rows.Item(144).RowHeight = 41.25
$ws.Rows.Item(145).RowHeight = 68.25
$ws.Rows.Item(146).RowHeight = 41.25
$ws.Rows.Item(147).RowHeight = 54.75
$ws.Rows.Item(148).RowHeight = 41.25
$ws.Rows.Item(149).RowHeight = 41.25
$ws.Rows.Item(150).RowHeight = 41.25
$ws.Rows.Item(151).RowHeight = 41.25
$ws.Rows.Item(152).RowHeight = 41.25

# ---------------------------------------------------------------------
# 4. Scroll / selection update so the freshly added rows are in view,
#    mirroring the sheetView topLeftCell/selection change in the diff.
# ---------------------------------------------------------------------

[void]$ws.Range("E155").Select()
$excel.ActiveWindow.ScrollRow = 149
$excel.ActiveWindow.ScrollColumn = 3
